$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.949892520904541
$ws.Range("B1").Value = 1.611510396003723
$ws.Range("C1").Value = 6.672751903533936
$ws.Range("D1").Value = 2.731826543807983
$ws.Range("E1").Value = 1.513468265533447
